$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1234603333333333
$ws.Range("H2").Value = 0.370381
$ws.Range("I2").Value = 0.002558470358543426
$ws.Range("J2").Value = 0.002636284444771545
$ws.Range("M2").Value = 9.363528666666667
$ws.Range("N2").Value = 28.090586
$ws.Range("O2").Value = 0.04175743631338733
$ws.Range("P2").Value = 0.04324026421082073
$ws.Range("Q2").Value = 1.156024370362889
$ws.Range("R2").Value = 10.404219333266
$ws.Range("S2").Value = 0.0001068351630565664
$ws.Range("T2").Value = 0.0001139936359267984

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1234603333333333
$ws.Range("H3").Value = 0.370381
$ws.Range("I3").Value = 0.002558470358543426
$ws.Range("J3").Value = 0.002636284444771545
$ws.Range("O3").Value = 0.1749266505387075
$ws.Range("P3").Value = 0.1811383852696593
$ws.Range("Q3").Value = 4.84271758282889
$ws.Range("R3").Value = 43.58445824546001
$ws.Range("S3").Value = 0.0004475446503225675
$ws.Range("T3").Value = 0.0004775323074374379

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1234603333333333
$ws.Range("H4").Value = 0.370381
$ws.Range("I4").Value = 0.002558470358543426
$ws.Range("J4").Value = 0.002636284444771545
$ws.Range("M4").Value = 71.284935
$ws.Range("N4").Value = 213.854805
$ws.Range("O4").Value = 0.3179011075133629
$ws.Range("P4").Value = 0.3291899382573772
$ws.Range("Q4").Value = 8.800861836745002
$ws.Range("R4").Value = 79.207756530705
$ws.Range("S4").Value = 0.0008133405605210658
$ws.Range("T4").Value = 0.0008678383136032289

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1234603333333333
$ws.Range("H5").Value = 0.370381
$ws.Range("I5").Value = 0.002558470358543426
$ws.Range("J5").Value = 0.002636284444771545
$ws.Range("M5").Value = 23.0690325
$ws.Range("N5").Value = 46.138065
$ws.Range("O5").Value = 0.1028782726814826
$ws.Range("P5").Value = 0.07102102180339065
$ws.Range("Q5").Value = 2.8481104421275
$ws.Range("R5").Value = 17.088662652765
$ws.Range("S5").Value = 0.0002632110111937211
$ws.Range("T5").Value = 0.0001872316150320595

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1234603333333333
$ws.Range("H6").Value = 0.370381
$ws.Range("I6").Value = 0.002558470358543426
$ws.Range("J6").Value = 0.002636284444771545
$ws.Range("M6").Value = 81.293813
$ws.Range("N6").Value = 243.881439
$ws.Range("O6").Value = 0.3625365329530597
$ws.Range("P6").Value = 0.3754103904587522
$ws.Range("Q6").Value = 10.03656125091767
$ws.Range("R6").Value = 90.329051258259
$ws.Range("S6").Value = 0.0009275389734495052
$ws.Range("T6").Value = 0.0009896885727720204

$ws.Range("I7").Value = 0.9088918061291337
$ws.Range("J7").Value = 0.9365351146153917
$ws.Range("M7").Value = 9.363528666666667
$ws.Range("N7").Value = 28.090586
$ws.Range("O7").Value = 0.04175743631338733
$ws.Range("P7").Value = 0.04324026421082073
$ws.Range("Q7").Value = 410.6754938160002
$ws.Range("R7").Value = 3696.079444344002
$ws.Range("S7").Value = 0.03795299171019689
$ws.Range("T7").Value = 0.04049602579868081

$ws.Range("I8").Value = 0.9088918061291337
$ws.Range("J8").Value = 0.9365351146153917
$ws.Range("O8").Value = 0.1749266505387075
$ws.Range("P8").Value = 0.1811383852696593
$ws.Range("S8").Value = 0.1589893993482456
$ws.Range("T8").Value = 0.1696424584097673

$ws.Range("I9").Value = 0.9088918061291337
$ws.Range("J9").Value = 0.9365351146153917
$ws.Range("M9").Value = 71.284935
$ws.Range("N9").Value = 213.854805
$ws.Range("O9").Value = 0.3179011075133629
$ws.Range("P9").Value = 0.3291899382573772
$ws.Range("Q9").Value = 3126.489694743265
$ws.Range("R9").Value = 28138.40725268938
$ws.Range("S9").Value = 0.2889377117782723
$ws.Range("T9").Value = 0.3082979365561065

$ws.Range("I10").Value = 0.9088918061291337
$ws.Range("J10").Value = 0.9365351146153917
$ws.Range("M10").Value = 23.0690325
$ws.Range("N10").Value = 46.138065
$ws.Range("O10").Value = 0.1028782726814826
$ws.Range("P10").Value = 0.07102102180339065
$ws.Range("Q10").Value = 1011.785903696867
$ws.Range("R10").Value = 6070.715422181204
$ws.Range("S10").Value = 0.09350521906891823
$ws.Range("T10").Value = 0.0665136807947407

$ws.Range("I11").Value = 0.9088918061291337
$ws.Range("J11").Value = 0.9365351146153917
$ws.Range("M11").Value = 81.293813
$ws.Range("N11").Value = 243.881439
$ws.Range("O11").Value = 0.3625365329530597
$ws.Range("P11").Value = 0.3754103904587522
$ws.Range("Q11").Value = 3565.469598743213
$ws.Range("R11").Value = 32089.22638868892
$ws.Range("S11").Value = 0.3295064842235007
$ws.Range("T11").Value = 0.3515850130560965

$ws.Range("G12").Value = 4.2730135
$ws.Range("H12").Value = 8.546027
$ws.Range("I12").Value = 0.08854972351232299
$ws.Range("J12").Value = 0.06082860093983664
$ws.Range("M12").Value = 9.363528666666667
$ws.Range("N12").Value = 28.090586
$ws.Range("O12").Value = 0.04175743631338733
$ws.Range("P12").Value = 0.04324026421082073
$ws.Range("Q12").Value = 40.01048440030367
$ws.Range("R12").Value = 240.062906401822
$ws.Range("S12").Value = 0.003697609440133884
$ws.Range("T12").Value = 0.002630244776213114

$ws.Range("G13").Value = 4.2730135
$ws.Range("H13").Value = 8.546027
$ws.Range("I13").Value = 0.08854972351232299
$ws.Range("J13").Value = 0.06082860093983664
$ws.Range("O13").Value = 0.1749266505387075
$ws.Range("P13").Value = 0.1811383852696593
$ws.Range("Q13").Value = 167.6084702626367
$ws.Range("R13").Value = 1005.65082157582
$ws.Range("S13").Value = 0.01548970654013929
$ws.Range("T13").Value = 0.01101839455245449

$ws.Range("G14").Value = 4.2730135
$ws.Range("H14").Value = 8.546027
$ws.Range("I14").Value = 0.08854972351232299
$ws.Range("J14").Value = 0.06082860093983664
$ws.Range("M14").Value = 71.284935
$ws.Range("N14").Value = 213.854805
$ws.Range("O14").Value = 0.3179011075133629
$ws.Range("P14").Value = 0.3291899382573772
$ws.Range("Q14").Value = 304.6014896016225
$ws.Range("R14").Value = 1827.608937609735
$ws.Range("S14").Value = 0.02815005517456955
$ws.Range("T14").Value = 0.02002416338766746

$ws.Range("G15").Value = 4.2730135
$ws.Range("H15").Value = 8.546027
$ws.Range("I15").Value = 0.08854972351232299
$ws.Range("J15").Value = 0.06082860093983664
$ws.Range("M15").Value = 23.0690325
$ws.Range("N15").Value = 46.138065
$ws.Range("O15").Value = 0.1028782726814826
$ws.Range("P15").Value = 0.07102102180339065
$ws.Range("Q15").Value = 98.57428730443876
$ws.Range("R15").Value = 394.297149217755
$ws.Range("S15").Value = 0.009109842601370655
$ws.Range("T15").Value = 0.004320109393617886

$ws.Range("G16").Value = 4.2730135
$ws.Range("H16").Value = 8.546027
$ws.Range("I16").Value = 0.08854972351232299
$ws.Range("J16").Value = 0.06082860093983664
$ws.Range("M16").Value = 81.293813
$ws.Range("N16").Value = 243.881439
$ws.Range("O16").Value = 0.3625365329530597
$ws.Range("P16").Value = 0.3754103904587522
$ws.Range("Q16").Value = 347.3695604154755
$ws.Range("R16").Value = 2084.217362492853
$ws.Range("S16").Value = 0.03210250975610961
$ws.Range("T16").Value = 0.02283568882988369
